$d = $word.ActiveDocument

# "Fix errore in 'inizio'": the closing sentence about DPR 151/11 was
# cut short; it should clarify that the non-subjection is "secondo
# quanto dichiarato" (according to what was declared), rather than
# ending flatly right after the decree reference.
#
# Locate the sentence-final "...DPR 151/11." and insert the new clause
# right before the trailing full stop, so the period that closes the
# paragraph stays its own run (matching how a real edit - placing the
# cursor just before the final period and typing the addition - would
# land in the document).
$findRng = $d.Content
$findRng.Find.Execute(
    "DPR 151/11.", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)

$periodStart = $findRng.End - 1
$insertionPoint = $d.Range($periodStart, $periodStart)
$insertionPoint.InsertBefore(", secondo quanto dichiarato")

# Nudge the newly inserted text onto its own run(s), distinct from the
# surrounding "... DPR 151/11" text and the trailing "." - mirrors the
# run boundaries Word leaves behind for freshly typed text.
$newTextRng = $d.Content
$newTextRng.Find.Execute(
    ", secondo quanto dichiarato", $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0
)
$newTextRng.Bold = $true
$newTextRng.Bold = $false
